# 4_particle_selection.xlsx
# Add an "actual size" column (with values for the two particles selected in
# weeks 1-2), add a "usable_tomer" review column, hide the rows that are not
# among the selected particles, and filter/mark the sheet by "usable".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before F ("video") for "actual size" ----------
$ws.Columns("F:F").Insert()

# --- 2. Insert a new (currently-beyond-range) column for "usable_tomer" ---
#        so the new cells inherit the existing cell style (s="1") the same
#        way the "actual size" column above did.
$ws.Columns("I:I").Insert()

# --- 3. Header row ----------------------------------------------------
$ws.Cells.Item(1, 6).Value = "actual size"
$ws.Cells.Item(1, 9).Value = "usable_tomer"

# --- 4. Column widths (C and D get wider) -------------------------------
$ws.Columns("C:C").ColumnWidth = 7.6
$ws.Columns("D:D").ColumnWidth = 6.6

# --- 5. "actual size" values for the two selected rows (weeks 1 & 2) ----
$ws.Cells.Item(2, 6).Value = 18.0
$ws.Cells.Item(12, 6).Value = 16.0

# --- 6. "usable" overrides (reviewer "tomer" disagreed on these rows) ----
$tomerOverrideRows = @(5, 6, 7, 8, 15, 16)
foreach ($r in $tomerOverrideRows) {
    $ws.Cells.Item($r, 8).Value = 0.0
}

# --- 7. "usable_tomer" flag for the rows tomer marked usable -------------
$usableTomerRows = @(6, 7, 8, 15, 16)
foreach ($r in $usableTomerRows) {
    $ws.Cells.Item($r, 9).Value = 1.0
}

# --- 8. Remove the leftover blank "usable_tomer" cells on every other row
$ws.Range("I2:I5").Clear()
$ws.Range("I9:I14").Clear()
$ws.Range("I17:I19").Clear()

# --- 9. Hide every data row except the two selected particles (2 and 12) -
for ($r = 3; $r -le 19; $r++) {
    if ($r -ne 12) {
        $ws.Rows("$r`:$r").Hidden = $true
    }
}

# --- 10. AutoFilter on the "usable" column (H), showing only val = 1 -----
$ws.Range("A1:H19").AutoFilter(8, @(1))

# --- 11. Hidden sheet-scoped defined name for the filter database --------
$name = $ws.Names.Add("_xlnm._FilterDatabase", "='4_particle_selection'!`$A`$1:`$H`$19")
$name.Visible = $false
